function Set-TextCell {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextCell $ws 2 4 "70.867.32"
Set-TextCell $ws 2 5 "  -0.19%  "

# Row 3
Set-TextCell $ws 3 4 "3.785.65"
Set-TextCell $ws 3 5 "  -1.20%  "

# Row 4
Set-TextCell $ws 4 5 "  +0.07%  "

# Row 5
Set-TextCell $ws 5 4 "697.23"
Set-TextCell $ws 5 5 "  -0.63%  "

# Row 6
Set-TextCell $ws 6 4 "169.19"
Set-TextCell $ws 6 5 "  -1.52%  "

# Row 7
Set-TextCell $ws 7 4 "3.786.10"
Set-TextCell $ws 7 5 "  -1.16%  "

# Row 8
Set-TextCell $ws 8 4 "1.01"
Set-TextCell $ws 8 5 "  +0.56%  "

# Row 9
Set-TextCell $ws 9 4 "0.521"
Set-TextCell $ws 9 5 "  -0.91%  "

# Row 10
Set-TextCell $ws 10 4 "0.159"
Set-TextCell $ws 10 5 "  -1.82%  "

# Row 11
Set-TextCell $ws 11 4 "7.53"
Set-TextCell $ws 11 5 "  +2.64%  "

# Row 12
Set-TextCell $ws 12 4 "0.476"
Set-TextCell $ws 12 5 "  +3.63%  "

# Row 13
Set-TextCell $ws 13 4 "0.0000248"
Set-TextCell $ws 13 5 "  -2.43%  "

# Row 14
Set-TextCell $ws 14 4 "35.98"
Set-TextCell $ws 14 5 "  -1.97%  "

# Row 15
Set-TextCell $ws 15 4 "4.436.33"
Set-TextCell $ws 15 5 "  -0.80%  "

# Row 16
Set-TextCell $ws 16 4 "3.882.46"
Set-TextCell $ws 16 5 "  +2.84%  "

# Row 17
Set-TextCell $ws 17 4 "71.143.46"
Set-TextCell $ws 17 5 "  +0.28%  "

# Row 18
Set-TextCell $ws 18 2 "TRON"
Set-TextCell $ws 18 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws 18 4 "0.115"
Set-TextCell $ws 18 5 "  -0.06%  "

# Row 19
Set-TextCell $ws 19 2 "Chainlink"
Set-TextCell $ws 19 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws 19 4 "17.48"
Set-TextCell $ws 19 5 "  +0.52%  "

# Row 20
Set-TextCell $ws 20 2 "Polkadot"
Set-TextCell $ws 20 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws 20 4 "7.15"
Set-TextCell $ws 20 5 "  -1.13%  "

# Row 21
Set-TextCell $ws 21 4 "514.70"
Set-TextCell $ws 21 5 "  +3.68%  "

# Row 22
Set-TextCell $ws 22 4 "10.36"
Set-TextCell $ws 22 5 "  -3.10%  "

# Row 23
Set-TextCell $ws 23 4 "0.711"
Set-TextCell $ws 23 5 "  -3.06%  "

# Row 24
Set-TextCell $ws 24 4 "83.47"
Set-TextCell $ws 24 5 "  -1.97%  "

# Row 25
Set-TextCell $ws 25 4 "0.0000139"
Set-TextCell $ws 25 5 "  -3.51%  "

# Row 26
Set-TextCell $ws 26 4 "12.52"
Set-TextCell $ws 26 5 "  +2.98%  "

# Row 27
Set-TextCell $ws 27 4 "3.943.79"
Set-TextCell $ws 27 5 "  -1.02%  "

# Row 28
Set-TextCell $ws 28 4 "10.17"
Set-TextCell $ws 28 5 "  -4.07%  "

# Row 29
Set-TextCell $ws 29 5 "  +0.04%  "

# Row 30
Set-TextCell $ws 30 4 "1.95"
Set-TextCell $ws 30 5 "  -6.58%  "

# Row 31
Set-TextCell $ws 31 4 "2.94"
Set-TextCell $ws 31 5 "  -4.06%  "

# Row 32
Set-TextCell $ws 32 4 "2.24"
Set-TextCell $ws 32 5 "  +0.43%  "

# Row 33
Set-TextCell $ws 33 4 "7.25"
Set-TextCell $ws 33 5 "  -2.72%  "

# Row 34
Set-TextCell $ws 34 4 "29.02"
Set-TextCell $ws 34 5 "  -1.25%  "

# Row 35
Set-TextCell $ws 35 5 "  -4.09%  "

# Row 36
Set-TextCell $ws 36 4 "9.22"
Set-TextCell $ws 36 5 "  +0.12%  "

# Row 37
Set-TextCell $ws 37 5 "  -1.12%  "

# Row 38
Set-TextCell $ws 38 4 "3.758.20"
Set-TextCell $ws 38 5 "  -0.90%  "

# Row 39
Set-TextCell $ws 39 4 "6.55"
Set-TextCell $ws 39 5 "  +9.73%  "

# Row 40
Set-TextCell $ws 40 4 "0.0998"
Set-TextCell $ws 40 5 "  -2.58%  "

# Row 41
Set-TextCell $ws 41 5 "  +1.34%  "

# Row 42
Set-TextCell $ws 42 5 "  -2.59%  "

# Row 44
Set-TextCell $ws 44 4 "3.17"
Set-TextCell $ws 44 5 "  -4.39%  "

# Row 45
Set-TextCell $ws 45 5 "  +0.25%  "

# Row 46
Set-TextCell $ws 46 4 "164.01"
Set-TextCell $ws 46 5 "  +0.48%  "

# Row 47
Set-TextCell $ws 47 4 "49.15"
Set-TextCell $ws 47 5 "  +0.32%  "

# Row 48
Set-TextCell $ws 48 4 "0.000300"
Set-TextCell $ws 48 5 "  -4.59%  "

# Row 49
Set-TextCell $ws 49 4 "418.06"
Set-TextCell $ws 49 5 "  -3.42%  "

# Row 50
Set-TextCell $ws 50 2 "Cosmos"
Set-TextCell $ws 50 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws 50 4 "8.59"
Set-TextCell $ws 50 5 "  -1.49%  "

# Row 51
Set-TextCell $ws 51 2 "ONDO"
Set-TextCell $ws 51 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell $ws 51 4 "1.37"
Set-TextCell $ws 51 5 "  -0.86%  "
